# Generate Report for Handback
# - Overview: zh-cn / de-de status columns flip from "Ready for handoff" to
#   "Handed back: in sync with en-US"
# - zh-cn / de-de detail sheets: bump "Latest Handback DateTime" and clear the
#   stale "Error Detail" (handback is now in sync, so no error remains)
# - A handful of column widths were widened/narrowed to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns for zh-cn (E2) and de-de (F2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet: Status (C2) shares the same string as Overview, plus
#     Latest Handback DateTime (K2) and Error Detail (P2) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-20 14:52:12"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- de-de sheet: Status (C2), Latest Handback DateTime (K2), Error Detail (P2) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-20 14:52:18"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334

